# Weekly update: a new price record for "Perejil" (Región de Ñuble / Terminal
# Hortofrutícola Agro Chillán) is inserted at the top of the data table
# (row 49, right below the header block that is shared across the other
# "Hortaliza" price rows already present in the sheet). Every existing
# data row from 49 down to 116 shifts down by one row (to 50..117); the
# previously-last row (old row 116) becomes the new last row, 117.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push all existing data rows (49:116) down by one to make room for the
# new weekly record.
$ws.Rows("49:49").Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A49").Value = 7
$ws.Range("B49").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C49").Value = "Ñuble"
$ws.Range("D49").Value = 45175
$ws.Range("E49").Value = 16
$ws.Range("F49").Value = 100112044
$ws.Range("G49").Value = "Perejil"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 180
$ws.Range("K49").Value = 1000
$ws.Range("L49").Value = 1000
$ws.Range("M49").Value = 1000
$ws.Range("N49").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O49").Value = "Región de Ñuble"
$ws.Range("P49").Value = 1000
$ws.Range("Q49").Value = 1
$ws.Range("R49").Value = "Hortaliza"
